$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4062.5
$ws.Range("I40").Value = 3800
$ws.Range("K40").Value = 3800
$ws.Range("M40").Value = -3625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 3959.7693
$ws.Range("I58").Value = 354
$ws.Range("J58").Value = 8166.5
$ws.Range("K58").Value = 1062
$ws.Range("L58").Value = 24499.5
$ws.Range("M58").Value = -912
$ws.Range("N58").Value = -24799.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3039.5833
$ws.Range("J80").Value = 3421.875
$ws.Range("L80").Value = 10265.625
$ws.Range("N80").Value = -12261.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 3039.5833
$ws.Range("J83").Value = 3421.875
$ws.Range("L83").Value = 30796.875
$ws.Range("N83").Value = -40780.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 440.1
$ws.Range("I101").Value = 356.1111
$ws.Range("J101").Value = 1196
$ws.Range("K101").Value = 1068.3333
$ws.Range("L101").Value = 3588
$ws.Range("M101").Value = 553.6667
$ws.Range("N101").Value = -6832

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 37580
$ws.Range("I131").Value = 2633.3333
$ws.Range("K131").Value = 7899.999899999999
$ws.Range("M131").Value = -2859.999899999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5594.8955
$ws.Range("I132").Value = 5777.4067
$ws.Range("J132").Value = 4248.875
$ws.Range("K132").Value = 17332.2201
$ws.Range("L132").Value = 12746.625
$ws.Range("M132").Value = -14802.2201
$ws.Range("N132").Value = -17806.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 8082.5
$ws.Range("J138").Value = 6200
$ws.Range("L138").Value = 18600
$ws.Range("N138").Value = -28880

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6854.4585
$ws.Range("I141").Value = 1282.0625
$ws.Range("J141").Value = 17999.25
$ws.Range("K141").Value = 3846.1875
$ws.Range("L141").Value = 53997.75
$ws.Range("M141").Value = 1333.8125
$ws.Range("N141").Value = -64357.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 204162.88
$ws.Range("I32").Value = 207646.14
$ws.Range("K32").Value = 207646.14
$ws.Range("M32").Value = -207359.14

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 33335594
$ws.Range("I122").Value = 83334500
$ws.Range("J122").Value = 2990.111
$ws.Range("K122").Value = 250003500
$ws.Range("L122").Value = 8970.332999999999
$ws.Range("M122").Value = -250001050
$ws.Range("N122").Value = -13870.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14215.6875
$ws.Range("I20").Value = 19994.6
$ws.Range("J20").Value = 4584.1665
$ws.Range("K20").Value = 19994.6
$ws.Range("L20").Value = 4584.1665
$ws.Range("M20").Value = -19747.6
$ws.Range("N20").Value = -5078.1665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1622.2307
$ws.Range("I94").Value = 751.9048
$ws.Range("K94").Value = 751.9048
$ws.Range("M94").Value = -300.9048

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7386.654
$ws.Range("I99").Value = 11016.1875
$ws.Range("K99").Value = 11016.1875
$ws.Range("M99").Value = -9518.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 12842962
$ws.Range("I134").Value = 25455.5
$ws.Range("K134").Value = 76366.5
$ws.Range("M134").Value = -73831.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 346876
$ws.Range("I16").Value = 409680.47
$ws.Range("J16").Value = 1451.5
$ws.Range("K16").Value = 409680.47
$ws.Range("L16").Value = 1451.5
$ws.Range("M16").Value = -409393.47
$ws.Range("N16").Value = -2025.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 346876
$ws.Range("I113").Value = 409680.47
$ws.Range("J113").Value = 1451.5
$ws.Range("K113").Value = 409680.47
$ws.Range("L113").Value = 1451.5
$ws.Range("M113").Value = -407510.47
$ws.Range("N113").Value = -5791.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 10110.083
$ws.Range("I122").Value = 2891.7144
$ws.Range("K122").Value = 8675.143199999999
$ws.Range("M122").Value = -6225.143199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2599
$ws.Range("I134").Value = 2559.1072
$ws.Range("J134").Value = 2822.4
$ws.Range("K134").Value = 7677.321599999999
$ws.Range("L134").Value = 8467.200000000001
$ws.Range("M134").Value = -5142.321599999999
$ws.Range("N134").Value = -13537.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 382406.12
$ws.Range("J141").Value = 422831.94
$ws.Range("L141").Value = 422831.94
$ws.Range("N141").Value = -433191.94

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I131").Value = 1510
$ws.Range("J131").Value = 7841.7856
$ws.Range("K131").Value = 4530
$ws.Range("L131").Value = 23525.3568
$ws.Range("M131").Value = 510
$ws.Range("N131").Value = -33605.3568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 8520.556
$ws.Range("I132").Value = 1439.6364
$ws.Range("J132").Value = 19647.715
$ws.Range("K132").Value = 12956.7276
$ws.Range("L132").Value = 176829.435
$ws.Range("M132").Value = -10426.7276
$ws.Range("N132").Value = -181889.435

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 4552.316
$ws.Range("I140").Value = 3540.9167
$ws.Range("J140").Value = 6286.143
$ws.Range("K140").Value = 10622.7501
$ws.Range("L140").Value = 18858.429
$ws.Range("M140").Value = -5442.750100000001
$ws.Range("N140").Value = -29218.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 6733629.5
$ws.Range("I2").Value = 6733629.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 6733629.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -6733516.5
$ws.Range("N2").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 83562.62
$ws.Range("I122").Value = 129543.875
$ws.Range("J122").Value = 9992.6
$ws.Range("K122").Value = 388631.625
$ws.Range("L122").Value = 29977.8
$ws.Range("M122").Value = -386181.625
$ws.Range("N122").Value = -34877.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3479
$ws.Range("I126").Value = 3479
$ws.Range("K126").Value = 10437
$ws.Range("M126").Value = -7967

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1012.6667
$ws.Range("I16").Value = 734.0714
$ws.Range("J16").Value = 1569.8572
$ws.Range("K16").Value = 734.0714
$ws.Range("L16").Value = 1569.8572
$ws.Range("M16").Value = -564.0714
$ws.Range("N16").Value = -1909.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 15533.533
$ws.Range("J61").Value = 1105
$ws.Range("L61").Value = 1105
$ws.Range("N61").Value = -1509

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 24980
$ws.Range("J64").Value = 24980
$ws.Range("L64").Value = 24980
$ws.Range("N64").Value = -25430

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 24980
$ws.Range("J67").Value = 24980
$ws.Range("L67").Value = 24980
$ws.Range("N67").Value = -26540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 15533.533
$ws.Range("J113").Value = 1105
$ws.Range("L113").Value = 1105
$ws.Range("N113").Value = -5445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 13167659
$ws.Range("I136").Value = 9620622
$ws.Range("K136").Value = 28861866
$ws.Range("M136").Value = -28859316

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 57999
$ws.Range("J63").Value = 57999
$ws.Range("L63").Value = 57999
$ws.Range("N63").Value = -59247

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 57999
$ws.Range("J66").Value = 57999
$ws.Range("L66").Value = 173997
$ws.Range("N66").Value = -180237

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 791.0909
$ws.Range("I100").Value = 585.6316
$ws.Range("J100").Value = 2092.3333
$ws.Range("K100").Value = 1171.2632
$ws.Range("L100").Value = 4184.6666
$ws.Range("M100").Value = -630.2632000000001
$ws.Range("N100").Value = -5266.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 5848902.5
$ws.Range("I113").Value = 761.3333
$ws.Range("K113").Value = 2283.9999
$ws.Range("M113").Value = -113.9998999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8397435
$ws.Range("I136").Value = 2072162.4
$ws.Range("K136").Value = 6216487.199999999
$ws.Range("M136").Value = -6213937.199999999
